$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $val) {
    $cell = $ws.Range($cellAddr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextValue "D2" '28.723.03'
Set-TextValue "E2" '  +2.32%  '
Set-TextValue "D3" '1.873.88'
Set-TextValue "E3" '  +2.24%  '
Set-TextValue "D4" '1.006'
Set-TextValue "E4" '  +0.44%  '
Set-TextValue "D5" '324.86'
Set-TextValue "E5" '  -0.31%  '
Set-TextValue "D6" '1.004'
Set-TextValue "E6" '  +0.32%  '
Set-TextValue "D7" '0.4594'
Set-TextValue "E7" '  -0.56%  '
Set-TextValue "D8" '0.3859'
Set-TextValue "E8" '  +0.02%  '
Set-TextValue "D9" '0.07862'
Set-TextValue "E9" '  +0.17%  '
Set-TextValue "D10" '0.9965'
Set-TextValue "E10" '  +3.71%  '
Set-TextValue "D11" '21.78'
Set-TextValue "E11" '  -0.76%  '
Set-TextValue "D12" '1.902.53'
Set-TextValue "E12" '  +0.66%  '
Set-TextValue "E13" '  +1.43%  '
Set-TextValue "D14" '5.699'
Set-TextValue "E14" '  +0.44%  '
Set-TextValue "D15" '0.06957'
Set-TextValue "E15" '  +1.69%  '
Set-TextValue "E16" '  +0.28%  '
Set-TextValue "E17" '  +0.42%  '
Set-TextValue "D18" '0.00001004'
Set-TextValue "E18" '  +1.15%  '
Set-TextValue "D19" '16.83'
Set-TextValue "E19" '  +0.84%  '
Set-TextValue "D20" '1.004'
Set-TextValue "E20" '  +0.24%  '
Set-TextValue "D21" '28.729.61'
Set-TextValue "E21" '  +2.23%  '
Set-TextValue "D22" '5.280'
Set-TextValue "E22" '  -0.19%  '
Set-TextValue "E23" '  +0.18%  '
Set-TextValue "D24" '2.125'
Set-TextValue "E24" '  +1.94%  '
Set-TextValue "D25" '2.115.66'
Set-TextValue "E25" '  +0.49%  '
Set-TextValue "D26" '153.59'
Set-TextValue "E26" '  -0.52%  '
Set-TextValue "D27" '19.24'
Set-TextValue "E27" '  +0.54%  '
Set-TextValue "D28" '5.762'
Set-TextValue "E28" '  +0.15%  '
Set-TextValue "E29" '  -0.35%  '
Set-TextValue "E30" '  +0.27%  '
Set-TextValue "D31" '0.09323'
Set-TextValue "E31" '  +0.85%  '
Set-TextValue "D32" '0.9187'
Set-TextValue "E32" '  -2.44%  '
Set-TextValue "D33" '5.309'
Set-TextValue "E33" '  +0.83%  '
Set-TextValue "D34" '1.338'
Set-TextValue "E34" '  +1.44%  '
Set-TextValue "D35" '3.311'
Set-TextValue "E35" '  -0.58%  '
Set-TextValue "D36" '0.05752'
Set-TextValue "E36" '  -1.16%  '
Set-TextValue "D37" '1.150'
Set-TextValue "E37" '  +1.32%  '
Set-TextValue "D38" '0.02073'
Set-TextValue "E38" '  -1.58%  '
Set-TextValue "D39" '7.721'
Set-TextValue "E39" '  -0.01%  '
Set-TextValue "D40" '0.5634'
Set-TextValue "E40" '  +0.73%  '
Set-TextValue "D41" '0.1789'
Set-TextValue "E41" '  +1.70%  '
Set-TextValue "D42" '9.897'
Set-TextValue "E42" '  -0.01%  '
Set-TextValue "D43" '0.07208'
Set-TextValue "E43" '  -1.58%  '
Set-TextValue "D44" '11.75'
Set-TextValue "E44" '  +1.19%  '
Set-TextValue "D45" '0.5281'
Set-TextValue "E45" '  +0.27%  '
Set-TextValue "D46" '2.134'
Set-TextValue "E46" '  +1.10%  '
Set-TextValue "D47" '1.117'
Set-TextValue "E47" '  -2.07%  '
Set-TextValue "D48" '113.51'
Set-TextValue "E48" '  +0.11%  '
Set-TextValue "D49" '1.825'
Set-TextValue "E49" '  -0.45%  '
Set-TextValue "E51" '  +0.38%  '
